$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1/IF1 ("I0" and "IF"), styled like the existing
# header row (bold, bordered, centered) by copying H1's formatting.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data values for the I/J columns in rows 2 and 3.
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8
